$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up row 4: remove empty inline-string cells E4:P4 ---
$ws.Range("E4:P4").ClearContents()

# --- Clean up row 5: remove empty inline-string cells E5:P5 ---
$ws.Range("E5:P5").ClearContents()

# --- Clean up row 6: remove empty B6 cell ---
$ws.Range("B6").ClearContents()

# --- Add row 7 ---
$ws.Range("A7").Value = "UVigoProfesor"
$ws.Range("B7").Value = "NO"
$ws.Range("C7").Value = "2025-10-09 16:30:19"
$ws.Range("D7").Value = "https://secretaria.uvigo.gal/uv/web/convocatoria/public/index"
$ws.Range("E7:P7").Value = "'"
$ws.Range("E7:P7").Style = "Normal"

# --- Add row 8 ---
$ws.Range("A8").Value = "USCEmprego"
$ws.Range("B8").Value = "ERROR"
$ws.Range("C8").Value = "2025-10-09 16:30:19"
$ws.Range("D8").Value = "https://www.usc.gal/gl/emprego"
$ws.Range("E8:P8").Value = "'"
$ws.Range("E8:P8").Style = "Normal"

# --- Add row 9 ---
$ws.Range("A9").Value = "USCConvocatorias"
$ws.Range("B9").Value = "'"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "2025-10-09 16:30:19"
$ws.Range("D9").Value = "https://www.usc.gal/gl/investigar-na-usc/convocatorias"
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 12
$ws.Range("H9").Value = 12
$ws.Range("I9").Value = 12
$ws.Range("J9").Value = 12
$ws.Range("K9").Value = 12
$ws.Range("L9").Value = 10
$ws.Range("M9").Value = 12
$ws.Range("N9").Value = 12
$ws.Range("O9").Value = 12
$ws.Range("P9").Value = 12
